# Update "想去人数" (F column) figures across sheets to the freshly
# generated values (gh-pages output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 46
$wsExhibit.Range("F3").Value = 19
$wsExhibit.Range("F4").Value = 206
$wsExhibit.Range("F5").Value = 4677
$wsExhibit.Range("F7").Value = 129
$wsExhibit.Range("F11").Value = 716
$wsExhibit.Range("F14").Value = 91
$wsExhibit.Range("F16").Value = 160
$wsExhibit.Range("F19").Value = 102
$wsExhibit.Range("F20").Value = 3684
$wsExhibit.Range("F21").Value = 6028
$wsExhibit.Range("F22").Value = 37
$wsExhibit.Range("F25").Value = 527
$wsExhibit.Range("F27").Value = 3406
$wsExhibit.Range("F28").Value = 378
$wsExhibit.Range("F30").Value = 2503
$wsExhibit.Range("F35").Value = 282
$wsExhibit.Range("F36").Value = 360
$wsExhibit.Range("F37").Value = 142
$wsExhibit.Range("F38").Value = 1538
$wsExhibit.Range("F41").Value = 38
$wsExhibit.Range("F44").Value = 474
$wsExhibit.Range("F45").Value = 70
$wsExhibit.Range("F46").Value = 557

# --- Sheet "演出" (Performances) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 100

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 46
$wsAll.Range("F3").Value = 19
$wsAll.Range("F4").Value = 206
$wsAll.Range("F5").Value = 4677
$wsAll.Range("F6").Value = 197
$wsAll.Range("F7").Value = 129
$wsAll.Range("F12").Value = 716
$wsAll.Range("F13").Value = 198
$wsAll.Range("F15").Value = 91
$wsAll.Range("F17").Value = 160
$wsAll.Range("F18").Value = 75
$wsAll.Range("F20").Value = 102
$wsAll.Range("F21").Value = 3684
$wsAll.Range("F22").Value = 6028
$wsAll.Range("F26").Value = 527
$wsAll.Range("F27").Value = 45
$wsAll.Range("F28").Value = 3406
$wsAll.Range("F29").Value = 378
$wsAll.Range("F31").Value = 2503
$wsAll.Range("F32").Value = 567
$wsAll.Range("F33").Value = 519
$wsAll.Range("F36").Value = 282
$wsAll.Range("F38").Value = 142
$wsAll.Range("F39").Value = 1538
$wsAll.Range("F41").Value = 27
$wsAll.Range("F42").Value = 38
$wsAll.Range("F43").Value = 53
$wsAll.Range("F44").Value = 465
$wsAll.Range("F45").Value = 474
$wsAll.Range("F47").Value = 557
